# Updates the "北京-漫展信息" workbook to the newly scraped data snapshot.
#
# Summary of the change:
#  - The exhibition "北京·thebONE✖️GOJO超次元动漫游戏嘉年华" (2024-01-20) has
#    dropped off the source feed, so its row is removed from both the
#    "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet. Every
#    row below it shifts up by one, and the running index in column A is
#    renumbered to stay sequential.
#  - The "想去人数" (interested-count) column F is refreshed with newer
#    counts across all four sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the "2024-01-20 thebONE GOJO" row from 展览 (row 3) and
#    全部类型 (row 5); everything below shifts up automatically.
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Rows(3).Delete()

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows(5).Delete()

# ---------------------------------------------------------------------
# 2. Renumber column A (the 0-based running index) so it stays
#    sequential after the deletion.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 31; $r++) {
    $wsExhibit.Cells.Item($r, 1).Value = ($r - 1)
}
for ($r = 2; $r -le 38; $r++) {
    $wsAll.Cells.Item($r, 1).Value = ($r - 1)
}

# ---------------------------------------------------------------------
# 3. Refresh the "想去人数" (column F) counts on every sheet.
# ---------------------------------------------------------------------

# 展览 (Exhibitions) - rows 2..31 after the deletion/shift above.
$exhibitF = @(30,93,830,510,271,9168,231,659,1709,35,55,2369,106,3628,259,95,110,180,222,180,77,40,240,473,98,1062,425,4272,56,211)
for ($i = 0; $i -lt $exhibitF.Length; $i++) {
    $wsExhibit.Cells.Item($i + 2, 6).Value = $exhibitF[$i]
}

# 演出 (Performances) - dimension unchanged, only F updates.
$wsPerf = $wb.Worksheets.Item("演出")
$wsPerf.Range("F2").Value = 36
$wsPerf.Range("F3").Value = 33
$wsPerf.Range("F5").Value = 13
$wsPerf.Range("F6").Value = 12

# 本地生活 (Local life) - dimension unchanged, only F updates.
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 179
$wsLocal.Range("F3").Value = 918

# 全部类型 (All types) - rows 2..38 after the deletion/shift above.
$allF = @(179,30,918,36,33,93,830,510,271,9168,231,659,1709,35,55,0,2369,106,3628,259,95,110,180,222,180,13,77,40,240,473,98,1062,425,4272,56,211,12)
for ($i = 0; $i -lt $allF.Length; $i++) {
    $wsAll.Cells.Item($i + 2, 6).Value = $allF[$i]
}
